$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (columns B-G) with new values ---
$data = @{
    2 = @(0.2156134072807172, 0.5104705812318155, 0.5597049646376884, 0.7481343225903276, 0.7243069924489576, 46)
    3 = @(0.07078259347676455, 0.3937967543584671, 0.2924574528699269, 0.5407933550534131, 0.5421993821439078, 45)
    4 = @(0.2022458033653479, 0.5029758090258937, 0.6072955914828451, 0.7792917242489138, 0.761291099874567, 44)
    5 = @(0.1056269426232762, 0.3853209831470272, 0.2947523084012032, 0.5429109580780288, 0.5388390533604958, 43)
    6 = @(0.2315116146747508, 0.5190355693975161, 0.6636775128766533, 0.8146640490881215, 0.7905440420027581, 42)
    7 = @(0.09265415981041521, 0.3465857693131973, 0.2074256174716431, 0.4554400262072308, 0.4514552522066898, 41)
    8 = @(0.2923186457898964, 0.5389143003675101, 0.734479721597224, 0.8570179237316008, 0.8158868265665942, 40)
    9 = @(0.1519280332023545, 0.335957238970132, 0.2134926424701896, 0.4620526403670794, 0.4420648348703681, 39)
    10 = @(0.2798477600321015, 0.4917781050680664, 0.5928569493473916, 0.7699720445233006, 0.7269447888063731, 38)
    11 = @(0.1565110834604237, 0.3372223941819412, 0.198541515633084, 0.4455799766967586, 0.4229425514680761, 37)
}

$labels = @{
    6 = "Q4"
    7 = "Q5"
    8 = "Q6"
    9 = "Q7"
    10 = "Q8"
    11 = "Q9"
}

# Update rows 2-5 in place
foreach ($r in 2..5) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}

# Add new rows 6-11, copying the label cell format from A5 then setting text,
# and filling in B-G numeric values.
$labelTemplate = $ws.Range("A5")
foreach ($r in 6..11) {
    $aCell = $ws.Cells.Item($r, 1)
    $labelTemplate.Copy($aCell)
    $aCell.Value = $labels[$r]

    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
